$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set updated / newly populated cell values
$ws.Range("C4").Value = 2015
$ws.Range("D4").Value = 2014
$ws.Range("E4").Value = 2013
$ws.Range("F4").Value = 2012
$ws.Range("C5").Value = 59.51
$ws.Range("D5").Value = 54.94
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("C6").Value = 1.37
$ws.Range("D6").Value = 2.49
$ws.Range("E6").Value = 27.23
$ws.Range("F6").Value = 24.97
$ws.Range("C7").Value = 1.96
$ws.Range("D7").Value = 1.43
$ws.Range("E7").Value = 4.24
$ws.Range("F7").Value = 4.76
$ws.Range("C8").Value = 37.16
$ws.Range("D8").Value = 41.14
$ws.Range("E8").Value = 68.53
$ws.Range("F8").Value = 70.27
$ws.Range("C16").Value = 2
$ws.Range("C17").Value = 2
$ws.Range("C18").Value = 4
$ws.Range("C28").Value = 0.38
$ws.Range("D28").Value = 0.63
$ws.Range("C29").Value = 0.5
$ws.Range("D29").Value = 0.5
$ws.Range("B38").Value = 2011
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 62.89
$ws.Range("B39").Value = 2012
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 43.39
$ws.Range("B40").Value = 2013
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 35.24
$ws.Range("B41").Value = 2014
$ws.Range("C41").Value = 6.25
$ws.Range("D41").Value = 42.76
$ws.Range("B42").Value = 2015
$ws.Range("C42").Value = 14.12
$ws.Range("D42").Value = 43.95
$ws.Range("C50").Value = 14.12
$ws.Range("D50").Value = "NA"
$ws.Range("C51").Value = 6.32
$ws.Range("D51").Value = 0.45
$ws.Range("B60").Value = 2013
$ws.Range("C60").Value = 0.1
$ws.Range("D60").Value = 1.39
$ws.Range("E60").Value = 8.36
$ws.Range("B61").Value = 2014
$ws.Range("C61").Value = 3.25
$ws.Range("D61").Value = 3.67
$ws.Range("E61").Value = 102.93
$ws.Range("B62").Value = 2015
$ws.Range("C62").Value = 4.1
$ws.Range("D62").Value = 6.5
$ws.Range("E62").Value = 73.31
$ws.Range("C82").Value = 2015
$ws.Range("D82").Value = 2014
$ws.Range("C83").Value = 0.84
$ws.Range("D83").Value = 0.84
$ws.Range("C84").Value = 0.64
$ws.Range("D84").Value = 0.64
$ws.Range("C85").Value = 0.05
$ws.Range("D85").Value = 0.05
$ws.Range("C96").Value = "FY 12/13"
$ws.Range("D96").Value = "FY 13/14"
$ws.Range("E96").Value = "FY 14/15"
$ws.Range("C97").Value = 0.47
$ws.Range("D97").Value = 0.84
$ws.Range("E97").Value = 6.29
$ws.Range("C98").Value = 0.14
$ws.Range("D98").Value = 0.64
$ws.Range("E98").Value = 6.44
$ws.Range("C99").Value = 2.11
$ws.Range("D99").Value = 0.05
$ws.Range("E99").Value = 0.3
$ws.Range("B110").Value = "FY 14/15"
$ws.Range("C110").Value = 14.12
$ws.Range("D110").Value = 43.95
$ws.Range("E110").Value = 1927.2
$ws.Range("B111").Value = "FY 13/14"
$ws.Range("C111").Value = 6.25
$ws.Range("D111").Value = 42.76
$ws.Range("E111").Value = 1076.09
$ws.Range("B112").Value = "FY 12/13"
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 35.24
$ws.Range("E112").Value = 120.77
$ws.Range("B113").Value = "FY 11/12"
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 43.39
$ws.Range("E113").Value = 1679.94
$ws.Range("B114").Value = "FY 10/11"
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 62.89
$ws.Range("E114").Value = 3432.8
$ws.Range("B115").Value = "FY 9/10"
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 47435
$ws.Range("E115").Value = 2188.09
$ws.Range("D123").Value = 0
$ws.Range("D124").Value = 41.25
$ws.Range("D125").Value = 0
$ws.Range("C134").Value = "FY 10/11"
$ws.Range("D134").Value = 0
$ws.Range("C135").Value = "FY 11/12"
$ws.Range("D135").Value = 0
$ws.Range("C136").Value = "FY 12/13"
$ws.Range("D136").Value = 0
$ws.Range("C137").Value = "FY 13/14"
$ws.Range("D137").Value = 0
$ws.Range("C138").Value = "FY 14/15"
$ws.Range("D138").Value = 41.25

# Clear cells whose values were removed
$ws.Range("C30").Value = $null
$ws.Range("C31").Value = $null
